# 31 de Marzo 2023
# Insert a new blank column before column A on the "ZAVALETA" sheet,
# shifting all existing data from A:H to B:I, then make that sheet
# the active/selected sheet.

$wb = $excel.ActiveWorkbook

$wsZavaleta = $wb.Worksheets.Item("ZAVALETA    ")

# Insert a new column before column A - shifts A:H data to B:I
$wsZavaleta.Columns.Item(1).Insert()

# Select the new active cell on the ZAVALETA sheet
$wsZavaleta.Range("F28").Select()

# Activate the ZAVALETA sheet (making it the active/selected tab)
$wsZavaleta.Activate()
